# Auto-generated Excel COM-interop edit script
# Applies the F-column (interest count) bumps across sheets 1-4,
# plus the sheet4 (全部类型) event-list content refresh for rows 19-21.

$wb = $excel.ActiveWorkbook

# --- Sheet 1 (展览) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 3285
$ws.Range("F4").Value = 16
$ws.Range("F5").Value = 714
$ws.Range("F6").Value = 2306
$ws.Range("F11").Value = 347
$ws.Range("F12").Value = 1093
$ws.Range("F13").Value = 453
$ws.Range("F15").Value = 87
$ws.Range("F16").Value = 248
$ws.Range("F17").Value = 4705
$ws.Range("F18").Value = 18
$ws.Range("F20").Value = 3506
$ws.Range("F22").Value = 126
$ws.Range("F24").Value = 3717
$ws.Range("F25").Value = 5120
$ws.Range("F27").Value = 979
$ws.Range("F29").Value = 3281
$ws.Range("F30").Value = 375
$ws.Range("F32").Value = 141
$ws.Range("F33").Value = 96
$ws.Range("F35").Value = 1197
$ws.Range("F37").Value = 14
$ws.Range("F38").Value = 1423
$ws.Range("F39").Value = 136
$ws.Range("F40").Value = 1381
$ws.Range("F41").Value = 883
$ws.Range("F45").Value = 339
$ws.Range("F46").Value = 74
$ws.Range("F49").Value = 3738

# --- Sheet 2 (演出) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F6").Value = 1015

# --- Sheet 3 (本地生活) ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 2273

# --- Sheet 4 (全部类型) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 2273
$ws.Range("F4").Value = 3285
$ws.Range("F5").Value = 16
$ws.Range("F6").Value = 714
$ws.Range("F8").Value = 2306
$ws.Range("F12").Value = 1015
$ws.Range("F15").Value = 347
$ws.Range("F16").Value = 1093
$ws.Range("F17").Value = 453
$ws.Range("F24").Value = 3506
$ws.Range("F25").Value = 3717
$ws.Range("F26").Value = 5120
$ws.Range("F28").Value = 979
$ws.Range("F30").Value = 375
$ws.Range("F32").Value = 141
$ws.Range("F33").Value = 96
$ws.Range("F35").Value = 1197
$ws.Range("F37").Value = 1423
$ws.Range("F38").Value = 136
$ws.Range("F39").Value = 1381
$ws.Range("F40").Value = 883
$ws.Range("F45").Value = 339
$ws.Range("F46").Value = 74
$ws.Range("F49").Value = 3738
$ws.Range("C19").Value = '杭州·重逢·怀旧only'
$ws.Range("D19").Value = '丰庆路492号建冠龙禾商务中心A幢 杭州华礼宴国际礼宴中心(龙禾商务中心店)'
$ws.Range("E19").Value = '2024.07.06 09:00-07.06 17:00'
$ws.Range("F19").Value = 87
$ws.Range("G19").Value = 69
$ws.Range("H19").Value = 'https://show.bilibili.com/platform/detail.html?id=85742'
$ws.Range("I19").Value = '//i2.hdslb.com/bfs/openplatform/202405/qBeP0pEz1715399357252.jpeg'
$ws.Range("B20").Value = "'2024-07-06"
$ws.Range("C20").Value = '杭州·黑执事only'
$ws.Range("D20").Value = '大岭山路156号 爱丽芬城堡'
$ws.Range("E20").Value = '2024.07.06 10:00-07.07 18:00'
$ws.Range("F20").Value = 248
$ws.Range("G20").Value = 160
$ws.Range("H20").Value = 'https://show.bilibili.com/platform/detail.html?id=86414'
$ws.Range("I20").Value = '//i1.hdslb.com/bfs/openplatform/202405/iP2cxk2w1716800288950.jpeg'
$ws.Range("C21").Value = '杭州·AD04动漫展'
$ws.Range("D21").Value = '阳城路雅澳杭州电商产业园西侧约200米 杭州大会展中心'
$ws.Range("E21").Value = '2024.07.13 10:00-07.14 17:00'
$ws.Range("F21").Value = 4705
$ws.Range("G21").Value = 75
$ws.Range("H21").Value = 'https://show.bilibili.com/platform/detail.html?id=85012'
$ws.Range("I21").Value = '//i0.hdslb.com/bfs/openplatform/202405/y1iKqqnh1715326769523.jpeg'
